$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory")
$ws.Select()

# Use an existing formatted row (row 8, columns C:N, style 63 = Text format + right align)
# as a template for the 5 new rows so the new cells inherit the same number format/style.
$ws.Range("C8:N8").Copy() | Out-Null
$ws.Range("C31:N35").PasteSpecial(-4122) | Out-Null
$ws.Range("C31:N35").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = $false

# Row 31 - new issue 1173 baseline memory reading
$ws.Range("A31").Value = 1173
$ws.Range("C31").Value = '$1FE4'
$ws.Range("D31").Value = '$7B00'
$ws.Range("E31").Value = '$BD00'
$ws.Range("F31").Value = "23324"
$ws.Range("G31").Value = '$199F'
$ws.Range("H31").Value = '$9A80'
$ws.Range("I31").Value = '$9A80'
$ws.Range("J31").Value = "32993"
$ws.Range("K31").Value = '$EC7F'
$ws.Range("L31").Value = '$EE2D'
$ws.Range("M31").Value = '$FFFA'
$ws.Range("N31").Value = "04557"

# Row 32 - netd
$ws.Range("B32").Value = "netd"
$ws.Range("C32").Value = '$1FE4'
$ws.Range("D32").Value = '$5500'
$ws.Range("E32").Value = '$BD00'
$ws.Range("F32").Value = "13596"
$ws.Range("G32").Value = '$199F'
$ws.Range("H32").Value = '$9A80'
$ws.Range("I32").Value = '$9A80'
$ws.Range("J32").Value = "32993"
$ws.Range("K32").Value = '$EC7F'
$ws.Range("L32").Value = '$F0D1'
$ws.Range("M32").Value = '$FFFA'
$ws.Range("N32").Value = "03881"

# Row 33 - drv uther
$ws.Range("B33").Value = "drv uther"
$ws.Range("C33").Value = '$1FE4'
$ws.Range("D33").Value = '$7B00'
$ws.Range("E33").Value = '$BD00'
$ws.Range("F33").Value = "23324"
$ws.Range("G33").Value = '$199F'
$ws.Range("H33").Value = '$9A80'
$ws.Range("I33").Value = '$9A80'
$ws.Range("J33").Value = "32993"
$ws.Range("K33").Value = '$EC7F'
$ws.Range("L33").Value = '$F0D1'
$ws.Range("M33").Value = '$FFFA'
$ws.Range("N33").Value = "03881"

# Row 34 - networkd
$ws.Range("B34").Value = "networkd"
$ws.Range("C34").Value = '$1FE4'
$ws.Range("D34").Value = '$4360'
$ws.Range("E34").Value = '$BD00'
$ws.Range("F34").Value = "09084"
$ws.Range("G34").Value = '$199F'
$ws.Range("H34").Value = '$9A80'
$ws.Range("I34").Value = '$9A80'
$ws.Range("J34").Value = "32993"
$ws.Range("K34").Value = '$EC7F'
$ws.Range("L34").Value = '$F0D1'
$ws.Range("M34").Value = '$FFFA'
$ws.Range("N34").Value = "03881"

# Row 35 - telnetd
$ws.Range("B35").Value = "telnetd"
$ws.Range("C35").Value = '$1FE4'
$ws.Range("D35").Value = '$40B0'
$ws.Range("E35").Value = '$BD00'
$ws.Range("F35").Value = "08396"
$ws.Range("G35").Value = '$199F'
$ws.Range("H35").Value = '$9A80'
$ws.Range("I35").Value = '$9A80'
$ws.Range("J35").Value = "32993"
$ws.Range("K35").Value = '$EC7F'
$ws.Range("L35").Value = '$F0D1'
$ws.Range("M35").Value = '$FFFA'
$ws.Range("N35").Value = "03881"

# New comment on B34 (networkd) matching the one already on B28
$cmt = $ws.Range("B34").AddComment("Patrick:" + [char]10 + "manually did networkd libtcpip &")

# Update selection to reflect the new active cell after the edit
$ws.Range("C36").Select() | Out-Null
